$d = $word.ActiveDocument

# Locate the paragraph block under the "DDD" bullet that currently reads
# (in order): Entities / Value objects - Immutable / Repositories /
# Indeholder database-access logien / Services / For logik der ikke falder
# entities naturligt.
# We find its first paragraph ("Entities") and its last paragraph ("For
# logik der ikke falder entities naturligt") by exact text match so we can
# replace that whole run of paragraphs in one shot - re-ordering the items
# and inserting the two new ones ("Stateless" and the split-off "Value
# objects" / "Immutable" bullets) as required by the change.
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($startPara -eq $null -and $t -eq "Entities ") { $startPara = $p }
    if ($t -eq "For logik der ikke falder entities naturligt") { $endPara = $p }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the DDD bullet block to update"
}

$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>
            </w:pPr>
            <w:r><w:t>Services</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr>
            </w:pPr>
            <w:r><w:t>For logik der ikke falder entities naturligt</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr>
              <w:rPr><w:b/></w:rPr>
            </w:pPr>
            <w:r><w:rPr><w:b/></w:rPr><w:t>Stateless</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>
            </w:pPr>
            <w:r><w:t>Repositories</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Entities </w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>
            </w:pPr>
            <w:r><w:t>Value objects</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr>
              <w:rPr><w:b/></w:rPr>
            </w:pPr>
            <w:r><w:rPr><w:b/></w:rPr><w:t>Immutable</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)

Write-Output "done"
